$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 500
$ws.Range("J13").Value = 500
$ws.Range("L13").Value = 500
$ws.Range("N13").Value = -838
$ws.Range("H17").Value = 1980.375
$ws.Range("I17").Value = 1850
$ws.Range("K17").Value = 5550
$ws.Range("M17").Value = -5382
$ws.Range("H28").Value = 730.9375
$ws.Range("I28").Value = 826.8182
$ws.Range("K28").Value = 826.8182
$ws.Range("M28").Value = -341.8182
$ws.Range("H31").Value = 316
$ws.Range("I31").Value = 316
$ws.Range("K31").Value = 948
$ws.Range("M31").Value = -718
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H125").Value = 2770.8572
$ws.Range("I125").Value = 2719.2
$ws.Range("K125").Value = 24472.8
$ws.Range("M125").Value = -22012.8
$ws.Range("H132").Value = 5349.1665
$ws.Range("I132").Value = 2136.8
$ws.Range("K132").Value = 6410.400000000001
$ws.Range("M132").Value = -3880.400000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4699.0713
$ws.Range("I32").Value = 3200
$ws.Range("K32").Value = 3200
$ws.Range("M32").Value = -2913
$ws.Range("H61").Value = 998
$ws.Range("I61").Value = 998
$ws.Range("K61").Value = 998
$ws.Range("M61").Value = -786
$ws.Range("H74").Value = 4125.6
$ws.Range("I74").Value = 4125.6
$ws.Range("K74").Value = 4125.6
$ws.Range("M74").Value = -3251.6
$ws.Range("H77").Value = 4125.6
$ws.Range("I77").Value = 4125.6
$ws.Range("K77").Value = 20628
$ws.Range("M77").Value = -16260
$ws.Range("H122").Value = 2272.0833
$ws.Range("I122").Value = 1973.8
$ws.Range("J122").Value = 2485.1428
$ws.Range("K122").Value = 5921.4
$ws.Range("L122").Value = 7455.428400000001
$ws.Range("M122").Value = -3471.4
$ws.Range("N122").Value = -12355.4284
$ws.Range("H136").Value = 998
$ws.Range("I136").Value = 998
$ws.Range("K136").Value = 2994
$ws.Range("M136").Value = -444

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 878.6667
$ws.Range("I64").Value = 954.4
$ws.Range("J64").Value = 500
$ws.Range("K64").Value = 954.4
$ws.Range("L64").Value = 500
$ws.Range("M64").Value = -729.4
$ws.Range("N64").Value = -950
$ws.Range("H67").Value = 878.6667
$ws.Range("I67").Value = 954.4
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 954.4
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = -174.4
$ws.Range("N67").Value = -2060
$ws.Range("H80").Value = 565.61536
$ws.Range("I80").Value = 467.77777
$ws.Range("J80").Value = 785.75
$ws.Range("K80").Value = 467.77777
$ws.Range("L80").Value = 785.75
$ws.Range("M80").Value = 530.2222300000001
$ws.Range("N80").Value = -2781.75
$ws.Range("H83").Value = 565.61536
$ws.Range("I83").Value = 467.77777
$ws.Range("J83").Value = 785.75
$ws.Range("K83").Value = 2338.88885
$ws.Range("L83").Value = 3928.75
$ws.Range("M83").Value = 2653.11115
$ws.Range("N83").Value = -13912.75
$ws.Range("H86").Value = 3753.8572
$ws.Range("I86").Value = 1944.25
$ws.Range("J86").Value = 6166.6665
$ws.Range("K86").Value = 1944.25
$ws.Range("L86").Value = 6166.6665
$ws.Range("M86").Value = -821.25
$ws.Range("N86").Value = -8412.666499999999
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 3753.8572
$ws.Range("I89").Value = 1944.25
$ws.Range("J89").Value = 6166.6665
$ws.Range("K89").Value = 9721.25
$ws.Range("L89").Value = 30833.3325
$ws.Range("M89").Value = -4105.25
$ws.Range("N89").Value = -42065.3325
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H134").Value = 999
$ws.Range("I134").Value = 999
$ws.Range("K134").Value = 2997
$ws.Range("M134").Value = -462

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 283.16666
$ws.Range("I7").Value = 259.8
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 259.8
$ws.Range("L7").Value = 400
$ws.Range("M7").Value = -146.8
$ws.Range("N7").Value = -626
$ws.Range("H17").Value = 994
$ws.Range("I17").Value = 994
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 994
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -820
$ws.Range("N17").ClearContents()
$ws.Range("H99").Value = 1366.6666
$ws.Range("I99").Value = 1300
$ws.Range("K99").Value = 1300
$ws.Range("M99").Value = 198
$ws.Range("H126").Value = 1366.6666
$ws.Range("I126").Value = 1300
$ws.Range("K126").Value = 3900
$ws.Range("M126").Value = -1430
$ws.Range("H132").Value = 6442.2
$ws.Range("I132").Value = 6427.75
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 19283.25
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -16753.25
$ws.Range("N132").Value = -24560
$ws.Range("H134").Value = 2164.4
$ws.Range("I134").Value = 2162.6667
$ws.Range("K134").Value = 6488.000100000001
$ws.Range("M134").Value = -3953.000100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3110
$ws.Range("I5").Value = 3800
$ws.Range("K5").Value = 11400
$ws.Range("M5").Value = -11288
$ws.Range("H14").Value = 911.2
$ws.Range("I14").Value = 911.2
$ws.Range("K14").Value = 2733.6
$ws.Range("M14").Value = -2560.6
$ws.Range("H68").Value = 1567.3334
$ws.Range("I68").Value = 1003
$ws.Range("J68").Value = 1849.5
$ws.Range("K68").Value = 3009
$ws.Range("L68").Value = 5548.5
$ws.Range("M68").Value = -2198
$ws.Range("N68").Value = -7170.5
$ws.Range("H71").Value = 1567.3334
$ws.Range("I71").Value = 1003
$ws.Range("J71").Value = 1849.5
$ws.Range("K71").Value = 9027
$ws.Range("L71").Value = 16645.5
$ws.Range("M71").Value = -4971
$ws.Range("N71").Value = -24757.5
$ws.Range("H107").Value = 622.5
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 663.3333
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 1989.9999
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -5829.9999
$ws.Range("H132").Value = 2500.5715
$ws.Range("I132").Value = 2500.8
$ws.Range("K132").Value = 22507.2
$ws.Range("M132").Value = -19977.2
$ws.Range("H135").Value = 3110
$ws.Range("I135").Value = 3800
$ws.Range("K135").Value = 34200
$ws.Range("M135").Value = -31665

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -10490
$ws.Range("H24").Value = 58006000
$ws.Range("I24").Value = 58006000
$ws.Range("K24").Value = 58006000
$ws.Range("M24").Value = -58005827
$ws.Range("H29").Value = 4500
$ws.Range("J29").Value = 4500
$ws.Range("L29").Value = 4500
$ws.Range("N29").Value = -5080
$ws.Range("H80").Value = 2476.625
$ws.Range("J80").Value = 2949.75
$ws.Range("L80").Value = 2949.75
$ws.Range("N80").Value = -4945.75
$ws.Range("H83").Value = 2476.625
$ws.Range("J83").Value = 2949.75
$ws.Range("L83").Value = 14748.75
$ws.Range("N83").Value = -24732.75
$ws.Range("H99").Value = 4900
$ws.Range("I99").Value = 4900
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4900
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2654
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 911.5
$ws.Range("I113").Value = 815.3333
$ws.Range("K113").Value = 815.3333
$ws.Range("M113").Value = 1354.6667
$ws.Range("H122").Value = 4592.3335
$ws.Range("I122").Value = 4910.8
$ws.Range("K122").Value = 14732.4
$ws.Range("M122").Value = -12282.4
$ws.Range("H132").Value = 3360.2856
$ws.Range("I132").Value = 3360.2856
$ws.Range("K132").Value = 10080.8568
$ws.Range("M132").Value = -7550.856800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 49995
$ws.Range("I74").Value = 49995
$ws.Range("K74").Value = 49995
$ws.Range("M74").Value = -48997
$ws.Range("H77").Value = 49995
$ws.Range("I77").Value = 49995
$ws.Range("K77").Value = 149985
$ws.Range("M77").Value = -144993

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 8000
$ws.Range("J22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("N22").Value = -8586
$ws.Range("H62").Value = 12500
$ws.Range("I62").Value = 13666.667
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 13666.667
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -13042.667
$ws.Range("N62").Value = -10248
$ws.Range("H65").Value = 12500
$ws.Range("I65").Value = 13666.667
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 68333.33499999999
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -65213.33499999999
$ws.Range("N65").Value = -51240
$ws.Range("H122").Value = 667966.7
$ws.Range("I122").Value = 667966.7
$ws.Range("K122").Value = 2003900.1
$ws.Range("M122").Value = -2001450.1
$ws.Range("H132").Value = 3852.5
$ws.Range("I132").Value = 3852.5
$ws.Range("K132").Value = 11557.5
$ws.Range("M132").Value = -9027.5
